$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.368.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "'1.872.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'235.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D7").Value = "'0.4669"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.2849"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").Value = "'21.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.39%  "
$ws.Range("D11").Value = "'0.07906"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").Value = "'98.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").Value = "'1.871.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "'5.130"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").Value = "'0.6762"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").Value = "'280.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "'30.369.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "'0.9998"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").Value = "'5.508"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.34%  "
$ws.Range("D20").Value = "'12.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").Value = "'2.110.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.21%  "
$ws.Range("D22").Value = "'0.000007309"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").Value = "'6.201"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("D25").Value = "'9.278"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").Value = "'165.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "'19.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("D28").Value = "'1.939"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.27%  "
$ws.Range("D29").Value = "'1.374"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "'0.09738"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").Value = "'4.431"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("D33").Value = "'4.122"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.75%  "
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("D35").Value = "'1.120"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.97%  "
$ws.Range("D36").Value = "'0.7072"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("D37").Value = "'2.714"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").Value = "'0.01862"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").Value = "'6.334"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.24%  "
$ws.Range("D40").Value = "'2.542"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("D41").Value = "'73.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("D42").Value = "'1.947"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("D43").Value = "'0.8505"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("D44").Value = "'0.4187"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "'104.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").Value = "'7.212"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").Value = "'9.163"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").Value = "'933.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.06%  "
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").Value = "'0.1132"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.21%  "
